# Fixed update to excel issue
#
# 1. Rename the "Requested quantity" header on the "Weekly Quantity" sheet
#    to "Weekly_PO_Qty".
# 2. Rename the "Requested quantity" header on the "Monthly Trend" sheet
#    to "Monthly_PO_Qty".
# 3. Add a new "PO Forecast" sheet (after "Monthly Trend") containing the
#    prophet-style forecast output: ds / PO_Forecast / yhat_lower / yhat_upper.

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# --- 1 & 2: header renames -------------------------------------------------
$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3: new "PO Forecast" sheet, placed after "Monthly Trend" -------------
$lastSheet  = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match the look & feel of the other sheets: bold/centered header row,
# date-number-formatted first column, same page margins & outline options.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2:A19").Copy()
$wsForecast.Range("A2:A27").PasteSpecial(-4122)

$wsForecast.Outline.SummaryRow    = 1
$wsForecast.Outline.SummaryColumn = 1

$wsForecast.PageSetup.LeftMargin   = 54
$wsForecast.PageSetup.RightMargin  = 54
$wsForecast.PageSetup.TopMargin    = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$forecastData = @(
    @(45347.99999999999, 417, -102.6513569269522, 903.1316906901793),
    @(45354.99999999999, 413, -102.5834085198147, 903.9389462467938),
    @(45375.99999999999, 402, -76.65664035817288, 889.9372722900764),
    @(45396.99999999999, 390, -107.4282057580699, 907.5728267546461),
    @(45410.99999999999, 383, -95.23222200254637, 871.4457600413206),
    @(45438.99999999999, 367, -131.010322688037, 820.4762827237058),
    @(45445.99999999999, 363, -103.9318267299962, 853.0573880738227),
    @(45480.99999999999, 344, -133.2161047396259, 859.8362344488233),
    @(45536.99999999999, 313, -165.1283852103562, 786.0654451202788),
    @(45543.99999999999, 310, -223.2092030224242, 798.071762077882),
    @(45550.99999999999, 306, -171.9326104120985, 807.2476896046303),
    @(45557.99999999999, 302, -176.0005177378167, 798.775907396898),
    @(45571.99999999999, 294, -207.1877559132853, 809.6922526781295),
    @(45578.99999999999, 290, -198.1399079107479, 789.4091392942778),
    @(45585.99999999999, 286, -229.5507129685772, 772.3795413645577),
    @(45592.99999999999, 283, -199.69484823218, 789.3808675344171),
    @(45599.99999999999, 279, -208.455847389776, 743.3020659938172),
    @(45606.99999999999, 275, -197.8689034359554, 781.4343836311368),
    @(45613.99999999999, 271, -205.7718692454265, 762.3217572668126),
    @(45620.99999999999, 267, -242.1863780630948, 707.9732706125274),
    @(45627.99999999999, 263, -243.3855166894361, 737.0714224098183),
    @(45634.99999999999, 260, -212.3948310429264, 755.5463683481086),
    @(45641.99999999999, 256, -248.3195200988411, 753.8409012259118),
    @(45648.99999999999, 252, -254.9959882282693, 710.5807866625681),
    @(45655.99999999999, 248, -268.5067697372983, 742.5677609709702),
    @(45662.99999999999, 244, -222.3260345720866, 727.9934233826754)
)

for ($i = 0; $i -lt $forecastData.Count; $i++) {
    $row    = $i + 2
    $values = $forecastData[$i]
    $wsForecast.Cells.Item($row, 1).Value = $values[0]
    $wsForecast.Cells.Item($row, 2).Value = $values[1]
    $wsForecast.Cells.Item($row, 3).Value = $values[2]
    $wsForecast.Cells.Item($row, 4).Value = $values[3]
}

# Leave the workbook's active sheet where it started (first sheet).
[void]$wsWeekly.Select()
